# "missed some parts for the BOM"
# The 1M resistor (row 12) and the 22k resistor (row 16) each actually need
# qty 2, not 1 - bump the Quantity column and let the Total column formulas
# (and the running TOTAL in K2) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12: 1M resistor quantity 1 -> 2
$ws.Range("F12").Value = 2

# Row 16: 22k resistor quantity 1 -> 2
$ws.Range("F16").Value = 2

# Re-apply the Total formula across the top block (H2:H8) so it consolidates
# into the same shared-formula group as the rest of the column (H9:H31),
# matching how Excel groups identical relative formulas together.
$ws.Range("H2:H8").Formula = "=(F2-G2)*E2"

# Leave the selection where it was left after editing row 25.
$ws.Range("I25").Select()
